$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Strain_Annotations")

# Delete the duplicate "Strain Descriptor" column (E) and shift everything left.
$ws.Range("E1").EntireColumn.Delete()

# Update header text for the remaining columns.
$ws.Range("B1").Value = "Strain Name(s)"
$ws.Range("H1").Value = "Parental Strain ID"

$ws.Range("C5").Select()
